$d = $word.ActiveDocument

# The document contains four "Warunki początkowe:" blocks; in each one we
# need to insert two new bulleted paragraphs right after the
# "Warunki początkowe:" heading paragraph (i.e. right before the existing
# "Użytkownik, który chce się zarejestrować nie może być zalogowany w
# systemie" bullet). Inserting a paragraph *before* that existing bullet
# paragraph makes Word clone its paragraph formatting (numbering/spacing/
# indent/run-size), which is exactly the formatting the new bullets need.

function Insert-BulletBefore($anchorIndex, $text) {
    $anchorRange = $d.Paragraphs.Item($anchorIndex).Range
    $anchorRange.InsertParagraphBefore()
    $d.Paragraphs.Item($anchorIndex).Range.Text = $text
}

function Add-StartingConditions($firstBulletIndex) {
    Insert-BulletBefore $firstBulletIndex "Użytkownik ma otwartą kartę aplikacji w przeglądarce"
    Insert-BulletBefore ($firstBulletIndex + 1) "Użytkownik znajduje się na stronie głównej aplikacji"
}

# Locate every paragraph whose text is the existing bullet that must now be
# preceded by the two new bullets, then process them back-to-front so
# earlier paragraph indices stay valid while we work.
$anchorText = "Użytkownik, który chce się zarejestrować nie może być zalogowany w systemie"
$anchorIndices = New-Object System.Collections.ArrayList
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a") -eq $anchorText) {
        [void]$anchorIndices.Add($i)
    }
}

if ($anchorIndices.Count -eq 0) {
    throw "Could not find any paragraph matching the expected anchor text."
}

for ($j = $anchorIndices.Count - 1; $j -ge 0; $j--) {
    Add-StartingConditions $anchorIndices[$j]
}
